$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1964.9286
$ws.Range("I38").Value = 114.5
$ws.Range("J38").Value = 4432.1665
$ws.Range("K38").Value = 343.5
$ws.Range("L38").Value = 13296.4995
$ws.Range("M38").Value = 28.5
$ws.Range("N38").Value = -14040.4995

$ws.Range("H74").Value = 8093.5386
$ws.Range("J74").Value = 8696.333000000001
$ws.Range("L74").Value = 8696.333000000001
$ws.Range("N74").Value = -10568.333

$ws.Range("H76").Value = 6269.75
$ws.Range("I76").Value = 6269.75
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6269.75
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5954.75
$ws.Range("N76").Value = $null

$ws.Range("H77").Value = 8093.5386
$ws.Range("J77").Value = 8696.333000000001
$ws.Range("L77").Value = 43481.665
$ws.Range("N77").Value = -52841.665

$ws.Range("H79").Value = 6269.75
$ws.Range("I79").Value = 6269.75
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6269.75
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5177.75
$ws.Range("N79").Value = $null

$ws.Range("H112").Value = 997.35
$ws.Range("J112").Value = 963.7222
$ws.Range("L112").Value = 2891.1666
$ws.Range("N112").Value = -5107.1666

$ws.Range("H132").Value = 3844
$ws.Range("I132").Value = 1574.9565
$ws.Range("J132").Value = 9642.666999999999
$ws.Range("K132").Value = 4724.8695
$ws.Range("L132").Value = 28928.001
$ws.Range("M132").Value = -2194.8695
$ws.Range("N132").Value = -33988.001

$ws.Range("H138").Value = 3604.2188
$ws.Range("I138").Value = 2747.75
$ws.Range("J138").Value = 3889.7083
$ws.Range("K138").Value = 8243.25
$ws.Range("L138").Value = 11669.1249
$ws.Range("M138").Value = -3103.25
$ws.Range("N138").Value = -21949.1249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2246.1755
$ws.Range("I32").Value = 1268.6735
$ws.Range("K32").Value = 1268.6735
$ws.Range("M32").Value = -981.6735000000001

$ws.Range("H61").Value = 2048.9722
$ws.Range("I61").Value = 1828.3549
$ws.Range("J61").Value = 3416.8
$ws.Range("K61").Value = 1828.3549
$ws.Range("L61").Value = 3416.8
$ws.Range("M61").Value = -1616.3549
$ws.Range("N61").Value = -3840.8

$ws.Range("H63").Value = 26800.8
$ws.Range("I63").Value = 8499.25
$ws.Range("J63").Value = 100007
$ws.Range("K63").Value = 8499.25
$ws.Range("L63").Value = 100007
$ws.Range("M63").Value = -7813.25
$ws.Range("N63").Value = -101379

$ws.Range("H66").Value = 26800.8
$ws.Range("I66").Value = 8499.25
$ws.Range("J66").Value = 100007
$ws.Range("K66").Value = 42496.25
$ws.Range("L66").Value = 500035
$ws.Range("M66").Value = -39064.25
$ws.Range("N66").Value = -506899

$ws.Range("H103").Value = 21475
$ws.Range("J103").Value = 21475
$ws.Range("L103").Value = 21475
$ws.Range("N103").Value = -23819

$ws.Range("H132").Value = 2235.2104
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 3761.2
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 11283.6
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -16343.6

$ws.Range("H136").Value = 2048.9722
$ws.Range("I136").Value = 1828.3549
$ws.Range("J136").Value = 3416.8
$ws.Range("K136").Value = 5485.0647
$ws.Range("L136").Value = 10250.4
$ws.Range("M136").Value = -2935.0647
$ws.Range("N136").Value = -15350.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1999.5
$ws.Range("I86").Value = 1999.5
$ws.Range("K86").Value = 1999.5
$ws.Range("M86").Value = -876.5

$ws.Range("H89").Value = 1999.5
$ws.Range("I89").Value = 1999.5
$ws.Range("K89").Value = 9997.5
$ws.Range("M89").Value = -4381.5

$ws.Range("H107").Value = 20838306
$ws.Range("I107").Value = 5775.077
$ws.Range("J107").Value = 111112610
$ws.Range("K107").Value = 5775.077
$ws.Range("L107").Value = 111112610
$ws.Range("M107").Value = -3855.077
$ws.Range("N107").Value = -111116450

$ws.Range("H134").Value = 3009.4644
$ws.Range("I134").Value = 1930.8928
$ws.Range("K134").Value = 5792.678400000001
$ws.Range("M134").Value = -3257.678400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6406.35
$ws.Range("I31").Value = 6433
$ws.Range("K31").Value = 6433
$ws.Range("M31").Value = -6138

$ws.Range("H34").Value = 6406.35
$ws.Range("I34").Value = 6433
$ws.Range("K34").Value = 6433
$ws.Range("M34").Value = -6231

$ws.Range("H134").Value = 2182.4614
$ws.Range("I134").Value = 1346.8214
$ws.Range("J134").Value = 4309.5454
$ws.Range("K134").Value = 4040.4642
$ws.Range("L134").Value = 12928.6362
$ws.Range("M134").Value = -1505.4642
$ws.Range("N134").Value = -17998.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45092190
$ws.Range("I4").Value = 50617680
$ws.Range("K4").Value = 151853040
$ws.Range("M4").Value = -151852928

$ws.Range("H129").Value = 1508.9231
$ws.Range("I129").Value = 731.44446
$ws.Range("J129").Value = 3258.25
$ws.Range("K129").Value = 2194.33338
$ws.Range("L129").Value = 9774.75
$ws.Range("M129").Value = 2805.66662
$ws.Range("N129").Value = -19774.75

$ws.Range("H131").Value = 1764.3636
$ws.Range("I131").Value = 695
$ws.Range("J131").Value = 1871.3
$ws.Range("K131").Value = 2085
$ws.Range("L131").Value = 5613.9
$ws.Range("M131").Value = 2955
$ws.Range("N131").Value = -15693.9

$ws.Range("H139").Value = 55184.633
$ws.Range("I139").Value = 61183.176
$ws.Range("K139").Value = 183549.528
$ws.Range("M139").Value = -178409.528

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11382.733
$ws.Range("J43").Value = 49998
$ws.Range("L43").Value = 49998
$ws.Range("N43").Value = -50300

$ws.Range("H80").Value = 603408.4
$ws.Range("I80").Value = 836716.3
$ws.Range("J80").Value = 253446.5
$ws.Range("K80").Value = 836716.3
$ws.Range("L80").Value = 253446.5
$ws.Range("M80").Value = -835718.3
$ws.Range("N80").Value = -255442.5

$ws.Range("H83").Value = 603408.4
$ws.Range("I83").Value = 836716.3
$ws.Range("J83").Value = 253446.5
$ws.Range("K83").Value = 4183581.5
$ws.Range("L83").Value = 1267232.5
$ws.Range("M83").Value = -4178589.5
$ws.Range("N83").Value = -1277216.5

$ws.Range("H102").Value = 2538.875
$ws.Range("I102").Value = 2518.8262
$ws.Range("K102").Value = 2518.8262
$ws.Range("M102").Value = -896.8262

$ws.Range("H126").Value = 2523.9583
$ws.Range("I126").Value = 2393.05
$ws.Range("J126").Value = 3178.5
$ws.Range("K126").Value = 7179.150000000001
$ws.Range("L126").Value = 9535.5
$ws.Range("M126").Value = -4709.150000000001
$ws.Range("N126").Value = -14475.5

$ws.Range("H132").Value = 23817702
$ws.Range("I132").Value = 25004760
$ws.Range("K132").Value = 75014280
$ws.Range("M132").Value = -75011750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900.46155
$ws.Range("I22").Value = 743.1667
$ws.Range("J22").Value = 1035.2858
$ws.Range("K22").Value = 743.1667
$ws.Range("L22").Value = 1035.2858
$ws.Range("M22").Value = -448.1667
$ws.Range("N22").Value = -1625.2858

$ws.Range("H27").Value = 900.46155
$ws.Range("I27").Value = 743.1667
$ws.Range("J27").Value = 1035.2858
$ws.Range("K27").Value = 743.1667
$ws.Range("L27").Value = 1035.2858
$ws.Range("M27").Value = -636.1667
$ws.Range("N27").Value = -1249.2858

$ws.Range("H46").Value = 983.3333
$ws.Range("I46").Value = 970
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 970
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -782
$ws.Range("N46").Value = -1376

$ws.Range("H68").Value = 3736.2144
$ws.Range("I68").Value = 2777.6667
$ws.Range("J68").Value = 5461.6
$ws.Range("K68").Value = 2777.6667
$ws.Range("L68").Value = 5461.6
$ws.Range("M68").Value = -2028.6667
$ws.Range("N68").Value = -6959.6

$ws.Range("H71").Value = 3736.2144
$ws.Range("I71").Value = 2777.6667
$ws.Range("J71").Value = 5461.6
$ws.Range("K71").Value = 13888.3335
$ws.Range("L71").Value = 27308
$ws.Range("M71").Value = -10144.3335
$ws.Range("N71").Value = -34796

$ws.Range("H101").Value = 62083.8
$ws.Range("J101").Value = 62083.8
$ws.Range("L101").Value = 62083.8
$ws.Range("N101").Value = -68573.8

$ws.Range("H127").Value = 148315.81
$ws.Range("J127").Value = 156147.4
$ws.Range("L127").Value = 156147.4
$ws.Range("N127").Value = -166067.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 86472
$ws.Range("J98").Value = 86472
$ws.Range("L98").Value = 86472
$ws.Range("N98").Value = -92462
